$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Link" URLs for rows that didn't have one yet ---
$ws.Range("D16").Value = "https://www.railwayoperationsimulator.com/wp-content/uploads/2019/06/Victoria-South-Eastern-scaled.jpg"
$ws.Range("D17").Value = "https://www.railwayoperationsimulator.com/wp-content/uploads/2018/02/Victoria-Central-Workstation-1.jpg"
$ws.Range("D18").Value = "https://www.railwayoperationsimulator.com/wp-content/uploads/2018/02/Victoria-Central-Workstation-2.jpg"

$ws.Range("D21").Value = "https://www.railwayoperationsimulator.com/wp-content/uploads/2018/02/Swanley-scaled.jpg"
$ws.Range("E21").Value = "Not including this project"

$ws.Range("D29").Value = "https://www.railwayoperationsimulator.com/wp-content/uploads/2018/07/Charing-Cross-Canon-St-London-Bridge-scaled.jpg"

# --- North Kent East Jn to Hayes and Dartford: status WIP -> TRUE (added) ---
$ws.Range("C41").Value = $true

# --- Grove Park, Hildenborough and Ashford: D68 becomes a real clickable hyperlink ---
$groveParkUrl = "https://www.railwayoperationsimulator.com/wp-content/uploads/2018/12/Grove-Pk-to-Bromley-North-and-Hildenborough.png"
$ws.Hyperlinks.Add($ws.Range("D68"), $groveParkUrl)
$ws.Range("D11").Copy()
$ws.Range("D68").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- South London and Thameslink: status WIP -> TRUE (added), plus new Link ---
$ws.Range("C86").Value = $true
$ws.Range("D86").Value = "https://www.railwayoperationsimulator.com/wp-content/uploads/2021/10/SouthLondonAndThameslink-scaled.bmp"

# --- Update view / selection state ---
$null = $ws.Range("C42").Select()
